$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.125.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.305.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.34"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.650"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0989"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.649.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.879"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.303.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.051.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  +6.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0848"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.58%  "
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "30.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.128"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0313"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "13.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.97%  "
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.26%  "
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.46%  "
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("E51").Value = "  -1.02%  "
